$wb = $excel.ActiveWorkbook

# --- Bonuses sheet: fill in Dutch translations for Name/Condition/Explanatory text/VP ---
$bonuses = $wb.Worksheets.Item("Bonuses")
$bonuses.Range("C2").Value = 'Ontleedkundige'
$bonuses.Range("D2").Value = 'Vogels met lichaamsdelen in hun naam'
$bonuses.Range("E2").Value = 'Onder lichaamsdelen vallen been, bek, borst, buik, eikel, hals, kaak, keel, kop, oog, oor, schouder, snavel, staart, teen, vleugel.'
$bonuses.Range("F2").Value = '2 of 3 vogels: 3[point]; 4+ vogels: 7[point]'
$bonuses.Range("C3").Value = 'Achtertuinvogelaar'
$bonuses.Range("D3").Value = 'Vogels die minder dan 4 punten waard zijn'
$bonuses.Range("F3").Value = '5 of 6 vogels: 3[point], 7+ vogels: 6[point]'
$bonuses.Range("C4").Value = 'Gedragswetenschapper'
$bonuses.Range("D4").Value = 'Voor elke kolom waarin vogels met eigenschappen in 3 verschillende kleuren liggen (bijvoorbeeld bruin, wit, roze, blauwgroen):'
$bonuses.Range("E4").Value = 'Vogels zonder eigenschap gelden als wit.'
$bonuses.Range("F4").Value = '3[point] per kolom'
$bonuses.Range("C5").Value = 'Vogelringer'
$bonuses.Range("D5").Value = 'Vogels die in meerdere leefomgevingen kunnen leven'
$bonuses.Range("F5").Value = '4 - 5 vogels: 4[point]; 6+ vogels: 7[point]'
$bonuses.Range("C6").Value = 'Vogelteller'
$bonuses.Range("D6").Value = 'Vogels met een [focking]-eigenschap'
$bonuses.Range("F6").Value = '2[point] per vogel'
$bonuses.Range("C7").Value = 'Vogelvoeder'
$bonuses.Range("D7").Value = 'Vogels die [seed] eten'
$bonuses.Range("E7").Value = 'Elke vogel met een [seed]-symbool. De vogel mag daarnaast andere soorten voedsel eten.'
$bonuses.Range("F7").Value = '5, 6 of 7 vogels: 3[point]; 8+ vogels: 7[point]'
$bonuses.Range("C8").Value = 'Broedmanager'
$bonuses.Range("D8").Value = 'Vogels waarop ten minste 4 eieren liggen'
$bonuses.Range("F8").Value = '1[point] per vogel'
$bonuses.Range("C9").Value = 'Kaartenmaker'
$bonuses.Range("D9").Value = 'Vogels met gebiedsgerelateerde termen in hun namen'
$bonuses.Range("E9").Value = 'Onder gebiedsgerelateerde namen vallen Amerikaans, Baltimore, berg, Californisch, Canadees, Carolina, Ĳslandse, Inca, Louisiana, Mexicaans, Mississippi, prairie, savannah, woestijn.'
$bonuses.Range("F9").Value = '3 of 4 vogels: 4[point]; 5+ vogels: 8[point]'
$bonuses.Range("C10").Value = 'Burgerwetenschapper'
$bonuses.Range("D10").Value = 'Vogels met weggestopte kaarten'
$bonuses.Range("F10").Value = '4 - 6 vogels: 3[point]; 7+ vogels: 6[point]'
$bonuses.Range("C11").Value = 'Dieetspecialist'
$bonuses.Range("D11").Value = 'Vogels met voedselkosten van 3 voedsel'
$bonuses.Range("F11").Value = '2 - 3 vogels: 3[point]; 4+ vogels: 6[point]'
$bonuses.Range("C12").Value = 'Ecoloog'
$bonuses.Range("D12").Value = 'Aantal vogels ini je leefomgeving met de minste vogels'
$bonuses.Range("E12").Value = 'Bij een gellijke stand: heb je bijvoorbeeld in elke leefomgeving 3 vogels, dan bevat je leefomgeving met de minste vogels er 3.'
$bonuses.Range("F12").Value = '2[point] per vogel'
$bonuses.Range("C13").Value = 'Gebiedsomheiner'
$bonuses.Range("D13").Value = 'Vogels met [ground]-nesten'
$bonuses.Range("E13").Value = 'De vogels moeten een [ground]- of [star]-nestsymbool hebben.'
$bonuses.Range("F13").Value = '4 of 5 vogels: 4[point]; 6+ vogels: 7[point]'
$bonuses.Range("C14").Value = 'Etholoog'
$bonuses.Range("D14").Value = 'In één leefomgevinig naar keuze:'
$bonuses.Range("E14").Value = '(bijvoobeeld bruin, wit, roze, blauwgroen)\nVogels zonder eigenschap gelden als wit.'
$bonuses.Range("F14").Value = '2[point] per eigenschapskleur'
$bonuses.Range("C15").Value = 'Valkenier'
$bonuses.Range("D15").Value = 'Vogels met een [predator]-eigenschap'
$bonuses.Range("F15").Value = '2[point] per vogel'
$bonuses.Range("C16").Value = 'Visserijbeheerder'
$bonuses.Range("D16").Value = 'Vogels die [fish] eten'
$bonuses.Range("E16").Value = 'Elke vogel met een [fish]-symbool. De vogel mag daarnaast andere soorten voedsel eten.'
$bonuses.Range("F16").Value = '2 of 3 vogels: 3[point]; 4+ vogels: 8[point]'
$bonuses.Range("C17").Value = 'Voedselwebexpert'
$bonuses.Range("D17").Value = 'Vogels die uitsluitend [invertebrate] eten'
$bonuses.Range("E17").Value = 'Elke vogel met een [invertebrate]-symbool en geen andere voedselsymbolen.'
$bonuses.Range("F17").Value = '2[point] per vogel'
$bonuses.Range("C18").Value = 'Boswachter'
$bonuses.Range("D18").Value = 'Vogels die uitsluitend in [forest] kunnen leven'
$bonuses.Range("F18").Value = '3 of 4 vogels: 4[point]; 5 vogels: 8[point]'
$bonuses.Range("C19").Value = 'Historicus'
$bonuses.Range("D19").Value = 'Vogels die naar een persoon zijn vernoemd'
$bonuses.Range("E19").Value = 'Onder deze categorie vallen Anna''s, Bairds, Bells, Bewicks, Brewers, Cassins, Clarks, Coopers, Forsters, Franklins, Gambels, Lincolns, Says, Stellers.'
$bonuses.Range("F19").Value = '2[point] per vogel'
$bonuses.Range("C20").Value = '"Grote Vogel"-Specialist'
$bonuses.Range("D20").Value = 'Vogels met vleugelwijdten boven 65 cm'
$bonuses.Range("F20").Value = '4 of 5 vogels: 3[point]; 6+ vogels: 6[point]'
$bonuses.Range("C21").Value = 'Nestkastenbouwer'
$bonuses.Range("D21").Value = 'Vogels met [cavity]-nesten'
$bonuses.Range("E21").Value = 'De vogels moeten een [cavity]- of [star]-nestsymbool hebben.'
$bonuses.Range("F21").Value = '4 of 5 vogels: 4[point]; 6+ vogels: 7[point]'
$bonuses.Range("C22").Value = 'Allesetersdeskundige'
$bonuses.Range("D22").Value = 'Vogels die [wild] eten'
$bonuses.Range("E22").Value = 'Elke vogel die een [wild]-symbool als onderdeel van zijn voedselkosten heeft.'
$bonuses.Range("F22").Value = '2[point] per vogel'
$bonuses.Range("C23").Value = 'Broedkundige'
$bonuses.Range("D23").Value = 'Vogels waar ten minste 1 ei op ligt'
$bonuses.Range("F23").Value = '7 of 8 vogels: 3[point]; 9+ vogels: 6[point]'
$bonuses.Range("C24").Value = 'Zangvogelspecalist'
$bonuses.Range("D24").Value = 'Vogels met vleugelwijdten van 30 cm of minder'
$bonuses.Range("F24").Value = '4 of 5 vogels: 3[point]; 6+ vogels: 6[point]'
$bonuses.Range("C25").Value = 'Fotograaf'
$bonuses.Range("D25").Value = 'Vogels met kleuren in hun namen'
$bonuses.Range("E25").Value = 'Onder kleuren vallen azuur, blauw, bont, bruin, citroen, geel, goud, grijs, groen, indigo, lazuli, purper, robijn, rood, rosse, tweekleurig, wit, zilver, zwart.'
$bonuses.Range("F25").Value = '4 of 5 vogels: 3[point]; 6+ vogels: 6[point]'
$bonuses.Range("C26").Value = 'Platformbouwer'
$bonuses.Range("D26").Value = 'Vogels met [platform]-nesten'
$bonuses.Range("E26").Value = 'De vogels moeten een [platform]- of [star]-nestsymbool hebben.'
$bonuses.Range("F26").Value = '4 of 5 vogels: 4[point]; 6+ vogels: 7[point]'
$bonuses.Range("C27").Value = 'Prairiebeheerder'
$bonuses.Range("D27").Value = 'Vogels die uitsluitend in [grassland] kunnen leven'
$bonuses.Range("F27").Value = '2 of 3 vogels: 3[point]; 4+ vogels: 8[point]'
$bonuses.Range("C28").Value = 'Knaagdierenexpert'
$bonuses.Range("D28").Value = 'Vogels die [rodent] eten'
$bonuses.Range("E28").Value = 'Elke vogel met een [rodent]-symbool. De vogel mag daarnaast andere soorten voedsel eten.'
$bonuses.Range("F28").Value = '2[point] per vogel'
$bonuses.Range("C29").Value = 'Visionair Leider'
$bonuses.Range("D29").Value = 'Vogelkaarten in de hand aan het einde van het spel'
$bonuses.Range("F29").Value = '5, 6 of 7 vogels: 4[point]; 8+ vogels: 7[point]'
$bonuses.Range("C30").Value = 'Wijnbouwer'
$bonuses.Range("D30").Value = 'Vogels die [fruit] eten'
$bonuses.Range("E30").Value = 'Elke vogel met een [fruit]-symbool. De vogel mag daarnaast andere soorten voedsel eten.'
$bonuses.Range("F30").Value = '2 of 3 vogels: 3[point]; 4+ vogels: 7[point]'
$bonuses.Range("C31").Value = 'Moeraswetenschapper'
$bonuses.Range("D31").Value = 'Vogels die uitsluitend in [wetland] kunnen leven'
$bonuses.Range("F31").Value = '3 of 4 vogels: 3[point]; 5 vogels: 7[point]'
$bonuses.Range("C32").Value = 'Heemtuinier'
$bonuses.Range("D32").Value = 'Vogels met [bowl]-nesten'
$bonuses.Range("E32").Value = 'De vogels moeten een [bowl]- of [star]-nestsymbool hebben.'
$bonuses.Range("F32").Value = '4 of 5 vogels: 4[point]; 6+ vogels: 7[point]'
$bonuses.Range("C33").Value = '[automa] Auvogelaar'
$bonuses.Range("D33").Value = 'Vogels die 3 of 4 punten waard zijn.'
$bonuses.Range("E33").Value = 'De Automa houdt er ten hoogste 2 (hogere waarde eerst).'
$bonuses.Range("C34").Value = '[automa] RASPB Erelid'
$bonuses.Range("D34").Value = 'Vogels die 5, 6 of 7 punten waard zijn.'
$bonuses.Range("E34").Value = 'De Automa houdt er ten hoogste één (die met de hoogste waarde).'

# --- Other sheet: fill in Dutch translations (Translated column) ---
$other = $wb.Worksheets.Item("Other")
$other.Range("B2").Value = 'ALS GEACTIVEERD'
$other.Range("B3").Value = 'ALS GESPEELD'
$other.Range("B4").Value = 'EENMALIG TUSSEN BEURTEN'
$other.Range("B5").Value = 'EINDE RONDE'
$other.Range("B6").Value = 'van de kaarten'

# --- View state: Bonuses tab becomes active/selected tab; update selections ---
$birds = $wb.Worksheets.Item("Birds")
$birds.Range("A160").Select()
$birds.Range("D185").Select()
$bonuses.Activate()
$bonuses.Range("E34").Select()
$other.Range("A7").Select()
